# edit.ps1 - applies the Time Tracking workbook update (#56 Time Tracking: updated to 2025-08-11)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1610
$ws.Range("A1610").Value = "2025-07-14"
$ws.Range("B1610").Value = "11:00"
$ws.Range("C1610").Value = "13:00"
$ws.Range("D1610").Value = "2h 00m"
$ws.Range("E1610").Value = "#maintenance"
$ws.Range("G1610").Value = "'False"
$ws.Range("H1610").Value = "'False"
$ws.Range("I1610").Formula = "=YEAR(A1610)"
$ws.Range("J1610").Formula = "=MONTH(A1610)"

# Row 1611
$ws.Range("A1611").Value = "2025-07-14"
$ws.Range("B1611").Value = "20:00"
$ws.Range("C1611").Value = "23:30"
$ws.Range("D1611").Value = "3h 30m"
$ws.Range("E1611").Value = "#maintenance"
$ws.Range("G1611").Value = "'False"
$ws.Range("H1611").Value = "'False"
$ws.Range("I1611").Formula = "=YEAR(A1611)"
$ws.Range("J1611").Formula = "=MONTH(A1611)"

# Row 1612
$ws.Range("A1612").Value = "2025-07-15"
$ws.Range("B1612").Value = "09:15"
$ws.Range("C1612").Value = "17:15"
$ws.Range("D1612").Value = "8h 00m"
$ws.Range("E1612").Value = "#python"
$ws.Range("F1612").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1612").Value = "'True"
$ws.Range("H1612").Value = "'False"
$ws.Range("I1612").Formula = "=YEAR(A1612)"
$ws.Range("J1612").Formula = "=MONTH(A1612)"

# Row 1613
$ws.Range("A1613").Value = "2025-07-16"
$ws.Range("B1613").Value = "21:15"
$ws.Range("C1613").Value = "22:45"
$ws.Range("D1613").Value = "1h 30m"
$ws.Range("E1613").Value = "#python"
$ws.Range("F1613").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1613").Value = "'True"
$ws.Range("H1613").Value = "'False"
$ws.Range("I1613").Formula = "=YEAR(A1613)"
$ws.Range("J1613").Formula = "=MONTH(A1613)"

# Row 1614
$ws.Range("A1614").Value = "2025-07-17"
$ws.Range("B1614").Value = "08:00"
$ws.Range("C1614").Value = "08:45"
$ws.Range("D1614").Value = "0h 45m"
$ws.Range("E1614").Value = "#studying"
$ws.Range("G1614").Value = "'False"
$ws.Range("H1614").Value = "'False"
$ws.Range("I1614").Formula = "=YEAR(A1614)"
$ws.Range("J1614").Formula = "=MONTH(A1614)"

# Row 1615
$ws.Range("A1615").Value = "2025-07-17"
$ws.Range("B1615").Value = "17:00"
$ws.Range("C1615").Value = "18:45"
$ws.Range("D1615").Value = "1h 45m"
$ws.Range("E1615").Value = "#studying"
$ws.Range("G1615").Value = "'False"
$ws.Range("H1615").Value = "'False"
$ws.Range("I1615").Formula = "=YEAR(A1615)"
$ws.Range("J1615").Formula = "=MONTH(A1615)"

# Row 1616
$ws.Range("A1616").Value = "2025-07-18"
$ws.Range("B1616").Value = "08:00"
$ws.Range("C1616").Value = "08:45"
$ws.Range("D1616").Value = "0h 45m"
$ws.Range("E1616").Value = "#studying"
$ws.Range("G1616").Value = "'False"
$ws.Range("H1616").Value = "'False"
$ws.Range("I1616").Formula = "=YEAR(A1616)"
$ws.Range("J1616").Formula = "=MONTH(A1616)"

# Row 1617
$ws.Range("A1617").Value = "2025-07-20"
$ws.Range("B1617").Value = "10:30"
$ws.Range("C1617").Value = "14:30"
$ws.Range("D1617").Value = "4h 00m"
$ws.Range("E1617").Value = "#python"
$ws.Range("F1617").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1617").Value = "'True"
$ws.Range("H1617").Value = "'False"
$ws.Range("I1617").Formula = "=YEAR(A1617)"
$ws.Range("J1617").Formula = "=MONTH(A1617)"

# Row 1618
$ws.Range("A1618").Value = "2025-07-20"
$ws.Range("B1618").Value = "15:30"
$ws.Range("C1618").Value = "17:30"
$ws.Range("D1618").Value = "2h 00m"
$ws.Range("E1618").Value = "#python"
$ws.Range("F1618").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1618").Value = "'True"
$ws.Range("H1618").Value = "'False"
$ws.Range("I1618").Formula = "=YEAR(A1618)"
$ws.Range("J1618").Formula = "=MONTH(A1618)"

# Row 1619
$ws.Range("A1619").Value = "2025-07-20"
$ws.Range("B1619").Value = "18:30"
$ws.Range("C1619").Value = "21:00"
$ws.Range("D1619").Value = "2h 30m"
$ws.Range("E1619").Value = "#python"
$ws.Range("F1619").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1619").Value = "'True"
$ws.Range("H1619").Value = "'False"
$ws.Range("I1619").Formula = "=YEAR(A1619)"
$ws.Range("J1619").Formula = "=MONTH(A1619)"

# Row 1620
$ws.Range("A1620").Value = "2025-07-20"
$ws.Range("B1620").Value = "21:30"
$ws.Range("C1620").Value = "23:00"
$ws.Range("D1620").Value = "1h 30m"
$ws.Range("E1620").Value = "#python"
$ws.Range("F1620").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1620").Value = "'True"
$ws.Range("H1620").Value = "'False"
$ws.Range("I1620").Formula = "=YEAR(A1620)"
$ws.Range("J1620").Formula = "=MONTH(A1620)"

# Row 1621
$ws.Range("A1621").Value = "2025-07-21"
$ws.Range("B1621").Value = "11:45"
$ws.Range("C1621").Value = "16:15"
$ws.Range("D1621").Value = "4h 30m"
$ws.Range("E1621").Value = "#python"
$ws.Range("F1621").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1621").Value = "'True"
$ws.Range("H1621").Value = "'False"
$ws.Range("I1621").Formula = "=YEAR(A1621)"
$ws.Range("J1621").Formula = "=MONTH(A1621)"

# Row 1622
$ws.Range("A1622").Value = "2025-07-21"
$ws.Range("B1622").Value = "21:15"
$ws.Range("C1622").Value = "00:00"
$ws.Range("D1622").Value = "2h 45m"
$ws.Range("E1622").Value = "#python"
$ws.Range("F1622").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1622").Value = "'True"
$ws.Range("H1622").Value = "'False"
$ws.Range("I1622").Formula = "=YEAR(A1622)"
$ws.Range("J1622").Formula = "=MONTH(A1622)"

# Row 1623
$ws.Range("A1623").Value = "2025-07-22"
$ws.Range("B1623").Value = "08:45"
$ws.Range("C1623").Value = "16:45"
$ws.Range("D1623").Value = "8h 00m"
$ws.Range("E1623").Value = "#python"
$ws.Range("F1623").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1623").Value = "'True"
$ws.Range("H1623").Value = "'False"
$ws.Range("I1623").Formula = "=YEAR(A1623)"
$ws.Range("J1623").Formula = "=MONTH(A1623)"

# Row 1624
$ws.Range("A1624").Value = "2025-07-22"
$ws.Range("B1624").Value = "20:00"
$ws.Range("C1624").Value = "22:00"
$ws.Range("D1624").Value = "2h 00m"
$ws.Range("E1624").Value = "#duckdb"
$ws.Range("F1624").Value = "nwapolloanalytics v1.0.0"
$ws.Range("G1624").Value = "'True"
$ws.Range("H1624").Value = "'False"
$ws.Range("I1624").Formula = "=YEAR(A1624)"
$ws.Range("J1624").Formula = "=MONTH(A1624)"

# Row 1625
$ws.Range("A1625").Value = "2025-07-24"
$ws.Range("B1625").Value = "08:00"
$ws.Range("C1625").Value = "08:45"
$ws.Range("D1625").Value = "0h 45m"
$ws.Range("E1625").Value = "#python"
$ws.Range("F1625").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1625").Value = "'True"
$ws.Range("H1625").Value = "'False"
$ws.Range("I1625").Formula = "=YEAR(A1625)"
$ws.Range("J1625").Formula = "=MONTH(A1625)"

# Row 1626
$ws.Range("A1626").Value = "2025-07-24"
$ws.Range("B1626").Value = "17:00"
$ws.Range("C1626").Value = "18:30"
$ws.Range("D1626").Value = "1h 30m"
$ws.Range("E1626").Value = "#python"
$ws.Range("F1626").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1626").Value = "'True"
$ws.Range("H1626").Value = "'False"
$ws.Range("I1626").Formula = "=YEAR(A1626)"
$ws.Range("J1626").Formula = "=MONTH(A1626)"

# Row 1627
$ws.Range("A1627").Value = "2025-07-25"
$ws.Range("B1627").Value = "08:00"
$ws.Range("C1627").Value = "08:45"
$ws.Range("D1627").Value = "0h 45m"
$ws.Range("E1627").Value = "#python"
$ws.Range("F1627").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1627").Value = "'True"
$ws.Range("H1627").Value = "'False"
$ws.Range("I1627").Formula = "=YEAR(A1627)"
$ws.Range("J1627").Formula = "=MONTH(A1627)"

# Row 1628
$ws.Range("A1628").Value = "2025-07-27"
$ws.Range("B1628").Value = "12:00"
$ws.Range("C1628").Value = "13:30"
$ws.Range("D1628").Value = "1h 30m"
$ws.Range("E1628").Value = "#maintenance"
$ws.Range("F1628").Value = "Resume update."
$ws.Range("G1628").Value = "'False"
$ws.Range("H1628").Value = "'False"
$ws.Range("I1628").Formula = "=YEAR(A1628)"
$ws.Range("J1628").Formula = "=MONTH(A1628)"

# Row 1629
$ws.Range("A1629").Value = "2025-07-28"
$ws.Range("B1629").Value = "10:00"
$ws.Range("C1629").Value = "15:00"
$ws.Range("D1629").Value = "5h 00m"
$ws.Range("E1629").Value = "#python"
$ws.Range("F1629").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1629").Value = "'True"
$ws.Range("H1629").Value = "'False"
$ws.Range("I1629").Formula = "=YEAR(A1629)"
$ws.Range("J1629").Formula = "=MONTH(A1629)"

# Row 1630
$ws.Range("A1630").Value = "2025-07-28"
$ws.Range("B1630").Value = "20:15"
$ws.Range("C1630").Value = "00:15"
$ws.Range("D1630").Value = "4h 00m"
$ws.Range("E1630").Value = "#python"
$ws.Range("F1630").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1630").Value = "'True"
$ws.Range("H1630").Value = "'False"
$ws.Range("I1630").Formula = "=YEAR(A1630)"
$ws.Range("J1630").Formula = "=MONTH(A1630)"

# Row 1631
$ws.Range("A1631").Value = "2025-07-29"
$ws.Range("B1631").Value = "10:15"
$ws.Range("C1631").Value = "12:45"
$ws.Range("D1631").Value = "2h 30m"
$ws.Range("E1631").Value = "#python"
$ws.Range("F1631").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1631").Value = "'True"
$ws.Range("H1631").Value = "'False"
$ws.Range("I1631").Formula = "=YEAR(A1631)"
$ws.Range("J1631").Formula = "=MONTH(A1631)"

# Row 1632
$ws.Range("A1632").Value = "2025-07-31"
$ws.Range("B1632").Value = "08:00"
$ws.Range("C1632").Value = "08:45"
$ws.Range("D1632").Value = "0h 45m"
$ws.Range("E1632").Value = "#maintenance"
$ws.Range("G1632").Value = "'False"
$ws.Range("H1632").Value = "'False"
$ws.Range("I1632").Formula = "=YEAR(A1632)"
$ws.Range("J1632").Formula = "=MONTH(A1632)"

# Row 1633
$ws.Range("A1633").Value = "2025-08-01"
$ws.Range("B1633").Value = "08:00"
$ws.Range("C1633").Value = "08:45"
$ws.Range("D1633").Value = "0h 45m"
$ws.Range("E1633").Value = "#python"
$ws.Range("F1633").Value = "nwrefurbishedanalytics v1.0.0"
$ws.Range("G1633").Value = "'True"
$ws.Range("H1633").Value = "'False"
$ws.Range("I1633").Formula = "=YEAR(A1633)"
$ws.Range("J1633").Formula = "=MONTH(A1633)"

# Row 1634
$ws.Range("A1634").Value = "2025-08-04"
$ws.Range("B1634").Value = "15:00"
$ws.Range("C1634").Value = "18:30"
$ws.Range("D1634").Value = "3h 30m"
$ws.Range("E1634").Value = "#maintenance"
$ws.Range("G1634").Value = "'False"
$ws.Range("H1634").Value = "'False"
$ws.Range("I1634").Formula = "=YEAR(A1634)"
$ws.Range("J1634").Formula = "=MONTH(A1634)"

# Row 1635
$ws.Range("A1635").Value = "2025-08-04"
$ws.Range("B1635").Value = "20:00"
$ws.Range("C1635").Value = "00:00"
$ws.Range("D1635").Value = "4h 00m"
$ws.Range("E1635").Value = "#maintenance"
$ws.Range("G1635").Value = "'False"
$ws.Range("H1635").Value = "'False"
$ws.Range("I1635").Formula = "=YEAR(A1635)"
$ws.Range("J1635").Formula = "=MONTH(A1635)"

# Row 1636
$ws.Range("A1636").Value = "2025-08-05"
$ws.Range("B1636").Value = "10:00"
$ws.Range("C1636").Value = "17:00"
$ws.Range("D1636").Value = "7h 00m"
$ws.Range("E1636").Value = "#maintenance"
$ws.Range("F1636").Value = "Hackberry Pi stuff."
$ws.Range("G1636").Value = "'False"
$ws.Range("H1636").Value = "'False"
$ws.Range("I1636").Formula = "=YEAR(A1636)"
$ws.Range("J1636").Formula = "=MONTH(A1636)"

# Row 1637
$ws.Range("A1637").Value = "2025-08-06"
$ws.Range("B1637").Value = "20:30"
$ws.Range("C1637").Value = "23:15"
$ws.Range("D1637").Value = "2h 45m"
$ws.Range("E1637").Value = "#maintenance"
$ws.Range("G1637").Value = "'False"
$ws.Range("H1637").Value = "'False"
$ws.Range("I1637").Formula = "=YEAR(A1637)"
$ws.Range("J1637").Formula = "=MONTH(A1637)"

# Row 1638
$ws.Range("A1638").Value = "2025-08-07"
$ws.Range("B1638").Value = "08:00"
$ws.Range("C1638").Value = "08:45"
$ws.Range("D1638").Value = "0h 45m"
$ws.Range("E1638").Value = "#studying"
$ws.Range("G1638").Value = "'False"
$ws.Range("H1638").Value = "'False"
$ws.Range("I1638").Formula = "=YEAR(A1638)"
$ws.Range("J1638").Formula = "=MONTH(A1638)"

# Row 1639
$ws.Range("A1639").Value = "2025-08-07"
$ws.Range("B1639").Value = "16:45"
$ws.Range("C1639").Value = "17:30"
$ws.Range("D1639").Value = "0h 45m"
$ws.Range("E1639").Value = "#studying"
$ws.Range("G1639").Value = "'False"
$ws.Range("H1639").Value = "'False"
$ws.Range("I1639").Formula = "=YEAR(A1639)"
$ws.Range("J1639").Formula = "=MONTH(A1639)"

# Row 1640
$ws.Range("A1640").Value = "2025-08-07"
$ws.Range("B1640").Value = "21:00"
$ws.Range("C1640").Value = "22:30"
$ws.Range("D1640").Value = "1h 30m"
$ws.Range("E1640").Value = "#studying"
$ws.Range("G1640").Value = "'False"
$ws.Range("H1640").Value = "'False"
$ws.Range("I1640").Formula = "=YEAR(A1640)"
$ws.Range("J1640").Formula = "=MONTH(A1640)"

# Row 1641
$ws.Range("A1641").Value = "2025-08-08"
$ws.Range("B1641").Value = "08:00"
$ws.Range("C1641").Value = "08:45"
$ws.Range("D1641").Value = "0h 45m"
$ws.Range("E1641").Value = "#studying"
$ws.Range("G1641").Value = "'False"
$ws.Range("H1641").Value = "'False"
$ws.Range("I1641").Formula = "=YEAR(A1641)"
$ws.Range("J1641").Formula = "=MONTH(A1641)"

# Row 1642
$ws.Range("A1642").Value = "2025-08-09"
$ws.Range("B1642").Value = "20:00"
$ws.Range("C1642").Value = "00:00"
$ws.Range("D1642").Value = "4h 00m"
$ws.Range("E1642").Value = "#maintenance"
$ws.Range("G1642").Value = "'False"
$ws.Range("H1642").Value = "'False"
$ws.Range("I1642").Formula = "=YEAR(A1642)"
$ws.Range("J1642").Formula = "=MONTH(A1642)"

# Row 1643
$ws.Range("A1643").Value = "2025-08-11"
$ws.Range("B1643").Value = "10:00"
$ws.Range("C1643").Value = "17:30"
$ws.Range("D1643").Value = "7h 30m"
$ws.Range("E1643").Value = "#maintenance"
$ws.Range("F1643").Value = "Debian 13 preseed.cfg."
$ws.Range("G1643").Value = "'False"
$ws.Range("H1643").Value = "'False"
$ws.Range("I1643").Formula = "=YEAR(A1643)"
$ws.Range("J1643").Formula = "=MONTH(A1643)"

# Row 1644
$ws.Range("A1644").Value = "2025-08-11"
$ws.Range("B1644").Value = "19:45"
$ws.Range("C1644").Value = "22:45"
$ws.Range("D1644").Value = "3h 00m"
$ws.Range("E1644").Value = "#maintenance"
$ws.Range("F1644").Value = "Debian 13 preseed.cfg."
$ws.Range("G1644").Value = "'False"
$ws.Range("H1644").Value = "'False"
$ws.Range("I1644").Formula = "=YEAR(A1644)"
$ws.Range("J1644").Formula = "=MONTH(A1644)"

# Update the view state (frozen pane scroll position + active selection) to match
# the author's final cursor position after entering the new rows.
$ws.Range("F1647").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1624
